$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.745.93"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "1.742.40"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").Value = "331.46"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").Value = "0.3837"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "0.3350"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").Value = "45.36"
$ws.Range("E9").Value = "  -5.13%  "
$ws.Range("D10").Value = "1.099"
$ws.Range("E10").Value = "  -3.45%  "
$ws.Range("D11").Value = "0.07163"
$ws.Range("E11").Value = "  -3.22%  "
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "22.09"
$ws.Range("E13").Value = "  -2.54%  "
$ws.Range("D14").Value = "6.104"
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").Value = "1.748.12"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").Value = "6.967"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").Value = "0.00001048"
$ws.Range("E17").Value = "  -2.59%  "
$ws.Range("D18").Value = "0.06578"
$ws.Range("E18").Value = "  -1.35%  "
$ws.Range("D19").Value = "1.003"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "78.35"
$ws.Range("E20").Value = "  -4.60%  "
$ws.Range("D21").Value = "16.63"
$ws.Range("E21").Value = "  -4.10%  "
$ws.Range("E22").Value = "  -4.03%  "
$ws.Range("D23").Value = "27.709.69"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "11.49"
$ws.Range("E24").Value = "  -4.60%  "
$ws.Range("D25").Value = "2.400"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "154.36"
$ws.Range("E26").Value = "  +0.48%  "
$ws.Range("D27").Value = "19.62"
$ws.Range("E27").Value = "  -5.23%  "
$ws.Range("D28").Value = "2.253"
$ws.Range("E28").Value = "  -6.69%  "
$ws.Range("D29").Value = "1.941.77"
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("D30").Value = "1.270"
$ws.Range("E30").Value = "  -11.57%  "
$ws.Range("D31").Value = "128.86"
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("D32").Value = "4.028"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").Value = "5.736"
$ws.Range("E33").Value = "  -6.57%  "
$ws.Range("D34").Value = "0.08680"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "11.94"
$ws.Range("E35").Value = "  -6.33%  "
$ws.Range("D37").Value = "5.087"
$ws.Range("E37").Value = "  -4.23%  "
$ws.Range("D38").Value = "0.02257"
$ws.Range("E38").Value = "  -7.12%  "
$ws.Range("D39").Value = "0.6432"
$ws.Range("E39").Value = "  -6.11%  "
$ws.Range("D40").Value = "0.06036"
$ws.Range("E40").Value = "  -4.52%  "
$ws.Range("D41").Value = "0.2080"
$ws.Range("E41").Value = "  -4.70%  "
$ws.Range("E42").Value = "  -3.74%  "
$ws.Range("D43").Value = "1.004"
$ws.Range("E43").Value = "  +0.69%  "
$ws.Range("D44").Value = "7.939"
$ws.Range("E44").Value = "  -3.59%  "
$ws.Range("D45").Value = "13.55"
$ws.Range("E45").Value = "  -4.20%  "
$ws.Range("D46").Value = "3.799"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "0.5952"
$ws.Range("E47").Value = "  -5.04%  "
$ws.Range("D48").Value = "125.83"
$ws.Range("E48").Value = "  -4.22%  "
$ws.Range("D49").Value = "1.972"
$ws.Range("E49").Value = "  -5.21%  "
$ws.Range("D50").Value = "1.143"
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("D51").Value = "0.06918"
$ws.Range("E51").Value = "  -6.19%  "
